$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 40.32364397296318
$ws.Range("G2").Value = 39.68057028333582
$ws.Range("H2").Value = 40.98317458832413
$ws.Range("I2").Value = 0.000761595484674244
$ws.Range("J2").Value = 0.0007110947592888342
$ws.Range("K2").Value = 0.0008564110728550661
$ws.Range("L2").Value = 0.05766856263294638
$ws.Range("M2").Value = 0.05723560521202441
$ws.Range("N2").Value = 0.05811785917621674
$ws.Range("F3").Value = 0.00001480428201179558
$ws.Range("G3").Value = 0.000000005511847492902186
$ws.Range("H3").Value = 0.00004197381530217415
$ws.Range("I3").Value = 0.00001296847474127172
$ws.Range("J3").Value = 0.000000005121929250915991
$ws.Range("K3").Value = 0.00003665257351095212
$ws.Range("L3").Value = 0.00001524844910802202
$ws.Range("M3").Value = 0.000000005745657754696681
$ws.Range("N3").Value = 0.00004321671004705058
$ws.Range("F4").Value = 40.32365877724519
$ws.Range("G4").Value = 39.68057028884766
$ws.Range("H4").Value = 40.98321656213943
$ws.Range("I4").Value = 0.0007745639594155157
$ws.Range("J4").Value = 0.000711099881218085
$ws.Range("K4").Value = 0.0008930636463660185
$ws.Range("L4").Value = 0.05768381108205441
$ws.Range("M4").Value = 0.05723561095768216
$ws.Range("N4").Value = 0.0581610758862638
